# feature/batch run and update log
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert a new "Description" column after run_number (old B..I -> C..J)
# ---------------------------------------------------------------------------
$ws.Columns("B:B").Insert()

# Remove the hyperlink that used to live on (old) C2 - its anchor cell
# reference does not move together with the column insert.
$ws.Range("C2").Hyperlinks.Delete()

# ---------------------------------------------------------------------------
# 2. Header row
# ---------------------------------------------------------------------------
$ws.Range("B1").Value = "Description"

# ---------------------------------------------------------------------------
# 3. Row 2 (Run 1) - update values in the now-shifted columns
# ---------------------------------------------------------------------------
$ws.Range("B2").Value = "Opening"
$ws.Range("C2").Value = "Business%20Development%20and%20Thought%20Leadership/Valuation%20Model%20Development/demo/Run1/Assumptions/Assumptions.xlsx"
$ws.Range("D2").Value = "https://datalyactuarial.sharepoint.com/sites/DatalyActuarial/Shared%20Documents/Business%20Development%20and%20Thought%20Leadership/Valuation%20Model%20Development/demo/Run1/models"
$ws.Range("E2").Value = "Business%20Development%20and%20Thought%20Leadership/Valuation%20Model%20Development/demo/Run1/model%20point%20files"
$ws.Range("F2").Value = "Business%20Development%20and%20Thought%20Leadership/Valuation%20Model%20Development/demo/Run1/outputs"
$ws.Range("G2").Value = 30
$ws.Range("H2").Value = [DateTime]"2024-06-30"
$ws.Range("I2").Value = "Basic_Term_Model_v0.4wRPG"
$ws.Range("J2").Value = "Model Points.xlsx"

# ---------------------------------------------------------------------------
# 4. Row 3 (Run 2)
# ---------------------------------------------------------------------------
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "Projection period from 30 to 20"
$ws.Range("C3").Value = "Business%20Development%20and%20Thought%20Leadership/Valuation%20Model%20Development/demo/Run2/Assumptions/Assumptions.xlsx"
$ws.Range("D3").Value = "https://datalyactuarial.sharepoint.com/sites/DatalyActuarial/Shared%20Documents/Business%20Development%20and%20Thought%20Leadership/Valuation%20Model%20Development/demo/Run2/models"
$ws.Range("E3").Value = "Business%20Development%20and%20Thought%20Leadership/Valuation%20Model%20Development/demo/Run2/model%20point%20files"
$ws.Range("F3").Value = "Business%20Development%20and%20Thought%20Leadership/Valuation%20Model%20Development/demo/Run2/outputs"
$ws.Range("G3").Value = 20
$ws.Range("H3").Value = [DateTime]"2024-06-30"
$ws.Range("I3").Value = "Basic_Term_Model_v0.4wRPG"
$ws.Range("J3").Value = "Model Points.xlsx"

# ---------------------------------------------------------------------------
# 5. Row 4 (Run 3) - brand new row
# ---------------------------------------------------------------------------
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "valiadation date change to end of 2024"
$ws.Range("C4").Value = "Business%20Development%20and%20Thought%20Leadership/Valuation%20Model%20Development/demo/Run3/Assumptions/Assumptions.xlsx"
$ws.Range("D4").Value = "https://datalyactuarial.sharepoint.com/sites/DatalyActuarial/Shared%20Documents/Business%20Development%20and%20Thought%20Leadership/Valuation%20Model%20Development/demo/Run3/models"
$ws.Range("E4").Value = "Business%20Development%20and%20Thought%20Leadership/Valuation%20Model%20Development/demo/Run3/model%20point%20files"
$ws.Range("F4").Value = "Business%20Development%20and%20Thought%20Leadership/Valuation%20Model%20Development/demo/Run3/outputs"
$ws.Range("G4").Value = 20
$ws.Range("H4").Value = [DateTime]"2024-12-31"
$ws.Range("I4").Value = "Basic_Term_Model_v0.4wRPG"
$ws.Range("J4").Value = "Model Points.xlsx"

# ---------------------------------------------------------------------------
# 6. Hyperlinks on the models_url column (D) for all three runs
# ---------------------------------------------------------------------------
$ws.Hyperlinks.Add($ws.Range("D2"), "https://datalyactuarial.sharepoint.com/sites/DatalyActuarial/Shared%20Documents/Business%20Development%20and%20Thought%20Leadership/Valuation%20Model%20Development/demo/Run1/models")
$ws.Range("D2").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("D3"), "https://datalyactuarial.sharepoint.com/sites/DatalyActuarial/Shared%20Documents/Business%20Development%20and%20Thought%20Leadership/Valuation%20Model%20Development/demo/Run2/models")
$ws.Range("D3").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("D4"), "https://datalyactuarial.sharepoint.com/sites/DatalyActuarial/Shared%20Documents/Business%20Development%20and%20Thought%20Leadership/Valuation%20Model%20Development/demo/Run3/models")
$ws.Range("D4").Style = "Hyperlink"

# ---------------------------------------------------------------------------
# 7. Column widths (best-effort match of the autofit widths in the diff)
# ---------------------------------------------------------------------------
$ws.Columns("A:A").ColumnWidth = 10
$ws.Columns("B:B").ColumnWidth = 10
$ws.Columns("C:C").ColumnWidth = 109
$ws.Columns("D:D").ColumnWidth = 157.17
$ws.Columns("E:E").ColumnWidth = 104.17
$ws.Columns("F:F").ColumnWidth = 91
$ws.Columns("G:G").ColumnWidth = 14.5
$ws.Columns("H:H").ColumnWidth = 11.83
$ws.Columns("I:I").ColumnWidth = 25.33
$ws.Columns("J:J").ColumnWidth = 21.5

# ---------------------------------------------------------------------------
# 8. Selection / view
# ---------------------------------------------------------------------------
$ws.Range("A2").Select()

# ---------------------------------------------------------------------------
# 9. Window size (bookViews)
# ---------------------------------------------------------------------------
$excel.Width = 32280
$excel.Height = 22020
